$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new (empty-ish) row 4 with a couple of values, like the original edit.
$ws.Range("A4").Value = "Test3"
$ws.Range("B4").Value = "123456"

# Move the active selection to J4, matching the saved view state.
$ws.Range("J4").Select()
